# Work against the "UT_Avance" template workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# Drop the two unused, empty helper sheets (Hoja2, Hoja3) - only Hoja1 stays.
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Hoja2").Delete()
$wb.Worksheets.Item("Hoja3").Delete()

# Push all existing report content down by one row, leaving a new blank
# row 1 at the top of the sheet.
$ws.Rows.Item(1).Insert()

# The blank spacer row that used to sit at row 5 (now shifted to row 6 by
# the insert above) is removed outright, closing the gap again so the rest
# of the report (old rows 6-12) lands back on its original row numbers.
$ws.Rows.Item(6).Delete()

# Leave the selection where the template was last left before saving.
$ws.Range("D25").Select()
